# Update the "Generate Report for Handback" timestamps.
# These cells hold text values that look like dates (formatted with a custom
# date numFmt), so they are stored as shared-string text, not numeric dates.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 ("Latest HO Xliff Generate Date" for the first file) and
# de-de!H2 ("Correspond Handoff Datetime" for the same file/locale pairing)
# shared the same timestamp text before the edit, and still share it after.
$wsOverview.Range("G2").Value = "2016-08-26 21:03:25"
$wsDeDe.Range("H2").Value     = "2016-08-26 21:03:25"

# zh-cn!H2 ("Correspond Handoff Datetime") and zh-cn!K2 ("Correspond Handback
# DateTime") for the first file.
$wsZhCn.Range("H2").Value = "2016-08-26 21:03:20"
$wsZhCn.Range("K2").Value = "2016-08-26 21:03:44"

# de-de!K2 ("Correspond Handback DateTime") for the first file.
$wsDeDe.Range("K2").Value = "2016-08-26 21:03:50"
